$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.Goto($ws.Range("D1"), $true)
